# Auto-applied update: BRVM - MAJ automatique via GitHub Actions
# Updates "Recommandations" and "Top_YTD" sheets with refreshed row
# values/order and appends 2 new rows (PALM CI, NESTLE CI) on Recommandations.

$wb = $excel.ActiveWorkbook

# --- Sheet: Recommandations ---
$ws1 = $wb.Worksheets.Item("Recommandations")

# Row 2
$ws1.Range("C2").Value = 8
$ws1.Range("D2").Value = 3388.72
$ws1.Range("E2").Value = 113.5

# Row 3
$ws1.Range("C3").Value = 4
$ws1.Range("D3").Value = 3280
$ws1.Range("E3").Value = 890

# Row 4
$ws1.Range("C4").Value = 4
$ws1.Range("D4").Value = 2795
$ws1.Range("E4").Value = 690

# Row 5
$ws1.Range("A5").Value = "BRVM - AUTRES SECTEURS"
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 2394.01
$ws1.Range("E5").Value = 590.09

# Row 6
$ws1.Range("A6").Value = "BRVM - DISTRIBUTION"
$ws1.Range("C6").Value = 4
$ws1.Range("D6").Value = 2365.61
$ws1.Range("E6").Value = 549.79

# Row 7
$ws1.Range("C7").Value = 4
$ws1.Range("D7").Value = 1452.73
$ws1.Range("E7").Value = 373.28

# Row 8
$ws1.Range("C8").Value = 4
$ws1.Range("D8").Value = 1345.8
$ws1.Range("E8").Value = 319.24

# Row 9
$ws1.Range("C9").Value = 4
$ws1.Range("D9").Value = 892.2
$ws1.Range("E9").Value = 202.58

# Row 10
$ws1.Range("A10").Value = "BRVM - FINANCES"
$ws1.Range("C10").Value = 4
$ws1.Range("D10").Value = 571.23
$ws1.Range("E10").Value = 141.89

# Row 11
$ws1.Range("A11").Value = "BRVM-PRESTIGE"
$ws1.Range("C11").Value = 4
$ws1.Range("D11").Value = 570.3200000000001
$ws1.Range("E11").Value = 140.27

# Row 12
$ws1.Range("C12").Value = 4
$ws1.Range("D12").Value = 561.4
$ws1.Range("E12").Value = 139.45

# Row 13
$ws1.Range("C13").Value = 4
$ws1.Range("D13").Value = 490.89
$ws1.Range("E13").Value = 120.52

# Row 14
$ws1.Range("C14").Value = 4
$ws1.Range("D14").Value = 427.06
$ws1.Range("E14").Value = 105.84

# Row 15
$ws1.Range("C15").Value = 4
$ws1.Range("D15").Value = 381.5
$ws1.Range("E15").Value = 95.17

# Row 19
$ws1.Range("A19").Value = "NEI-CEDA CI (NEIC)"
$ws1.Range("B19").Value = 4
$ws1.Range("D19").Value = 28.69
$ws1.Range("E19").Value = 7.34

# Row 20
$ws1.Range("A20").Value = "SOLIBRA CI (SLBC)"
$ws1.Range("D20").Value = 22.47
$ws1.Range("E20").Value = 7.49

# Row 22
$ws1.Range("A22").Value = "UNILEVER CI (UNLC)"
$ws1.Range("D22").Value = 6.98
$ws1.Range("E22").Value = 6.98

# Row 23
$ws1.Range("A23").Value = "SICOR CI (SICC)"
$ws1.Range("D23").Value = 6.06
$ws1.Range("E23").Value = 6.06

# Row 24
$ws1.Range("A24").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Range("D24").Value = 4.55
$ws1.Range("E24").Value = 4.55

# Row 25
$ws1.Range("A25").Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Range("B25").Value = 2
$ws1.Range("D25").Value = 3.71
$ws1.Range("E25").Value = 2.35

# Row 26
$ws1.Range("A26").Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Range("D26").Value = 2.59
$ws1.Range("E26").Value = 2.59

# Row 27
$ws1.Range("A27").Value = "SITAB CI (STBC)"
$ws1.Range("B27").Value = 1
$ws1.Range("C27").Value = 0
$ws1.Range("D27").Value = 2.56
$ws1.Range("E27").Value = 2.56

# Row 28
$ws1.Range("A28").Value = "TOTAL"
$ws1.Range("B28").Value = 0
$ws1.Range("C28").Value = 4
$ws1.Range("D28").Value = 0
$ws1.Range("E28").Value = 0
$ws1.Range("G28").Value = "➖ Neutre"

# Row 29
$ws1.Range("A29").Value = "CFAO MOTORS CI (CFAC)"
$ws1.Range("B29").Value = 1
$ws1.Range("D29").Value = -0.4
$ws1.Range("E29").Value = -7.48
$ws1.Range("G29").Value = "👀 À surveiller"

# Row 30
$ws1.Range("A30").Value = "SAPH CI (SPHC)"
$ws1.Range("B30").Value = 1
$ws1.Range("D30").Value = -2
$ws1.Range("E30").Value = 2.53
$ws1.Range("G30").Value = "👀 À surveiller"

# Row 31
$ws1.Range("A31").Value = "BANK OF AFRICA ML (BOAM)"
$ws1.Range("C31").Value = 1
$ws1.Range("D31").Value = -2.82
$ws1.Range("E31").Value = 4.65

# Row 32
$ws1.Range("A32").Value = "SOGB CI (SOGC)"
$ws1.Range("D32").Value = -3.83
$ws1.Range("E32").Value = -3.83

# Row 33
$ws1.Range("A33").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Range("D33").Value = -4.63
$ws1.Range("E33").Value = -4.63

# Row 34
$ws1.Range("A34").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Range("B34").Value = 1
$ws1.Range("C34").Value = 2
$ws1.Range("D34").Value = -7.43
$ws1.Range("E34").Value = -7.42
$ws1.Range("G34").Value = "👀 À surveiller"

# Row 35
$ws1.Range("A35").Value = "BERNABE CI (BNBC)"
$ws1.Range("D35").Value = -14.92
$ws1.Range("E35").Value = -7.48

# Row 36
$ws1.Range("A36").Value = "SAFCA CI (SAFC)"
$ws1.Range("E36").Value = -7.46

# Row 38
$ws1.Range("A38").Value = "PALM CI (PALC)"
$ws1.Range("B38").Value = 0
$ws1.Range("C38").Value = 3
$ws1.Range("D38").Value = -22.38
$ws1.Range("E38").Value = -7.44
$ws1.Range("F38").Value = "🔴 Vente"
$ws1.Range("G38").Value = "⚠️ Risque de décrochage"

# Row 39
$ws1.Range("A39").Value = "NESTLE CI (NTLC)"
$ws1.Range("B39").Value = 0
$ws1.Range("C39").Value = 3
$ws1.Range("D39").Value = -22.43
$ws1.Range("E39").Value = -7.47
$ws1.Range("F39").Value = "🔴 Vente"
$ws1.Range("G39").Value = "⚠️ Risque de décrochage"

# --- Sheet: Top_YTD ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

# Row 2
$ws2.Range("B2").Value = 10097505.15

# Row 3
$ws2.Range("B3").Value = 708297.72

# Row 4
$ws2.Range("B4").Value = 406876.4

# Row 5
$ws2.Range("A5").Value = "BRVM - AUTRES SECTEURS"
$ws2.Range("B5").Value = 237921.49

# Row 6
$ws2.Range("A6").Value = "BRVM - DISTRIBUTION"
$ws2.Range("B6").Value = 227459.44

# Row 7
$ws2.Range("B7").Value = 45906.91

# Row 8
$ws2.Range("B8").Value = 36135.33

# Row 9
$ws2.Range("B9").Value = 10740.49

# Row 10
$ws2.Range("A10").Value = "BRVM - FINANCES"
$ws2.Range("B10").Value = 3375.62

# Row 11
$ws2.Range("A11").Value = "BRVM-PRESTIGE"
$ws2.Range("B11").Value = 3362.37
